# Updates crypto prices/volumes to match the latest scrape (per commit diff).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    # Force text storage so numeric-looking strings (e.g. "539.13") keep their
    # exact formatting instead of being parsed into floating point numbers.
    $rng.NumberFormat = "@"
    $rng.Value = $value
}

# Row 2
$ws.Range("D2").Value = "58.977.93"
$ws.Range("E2").Value = "  -6.16%  "

# Row 3
$ws.Range("D3").Value = "2.443.99"
$ws.Range("E3").Value = "  -8.87%  "

# Row 4
$ws.Range("E4").Value = "  +0.01%  "

# Row 5
Set-TextCell "D5" "539.13"
$ws.Range("E5").Value = "  -2.67%  "

# Row 6
Set-TextCell "D6" "147.06"
$ws.Range("E6").Value = "  -7.12%  "

# Row 7
Set-TextCell "D7" "0.998"
$ws.Range("E7").Value = "  -0.15%  "

# Row 8
Set-TextCell "D8" "0.572"
$ws.Range("E8").Value = "  -3.22%  "

# Row 9
$ws.Range("D9").Value = "2.458.78"
$ws.Range("E9").Value = "  -8.43%  "

# Row 10
Set-TextCell "D10" "0.0991"
$ws.Range("E10").Value = "  -6.71%  "

# Row 11
$ws.Range("E11").Value = "  -2.21%  "

# Row 12
$ws.Range("E12").Value = "  -0.73%  "

# Row 13
Set-TextCell "D13" "0.351"
$ws.Range("E13").Value = "  -4.82%  "

# Row 14
$ws.Range("D14").Value = "2.882.91"
$ws.Range("E14").Value = "  -8.62%  "

# Row 15
Set-TextCell "D15" "23.93"
$ws.Range("E15").Value = "  -9.77%  "

# Row 16
$ws.Range("D16").Value = "58.882.17"
$ws.Range("E16").Value = "  -6.19%  "

# Row 17
$ws.Range("E17").Value = "  -6.44%  "

# Row 18
$ws.Range("D18").Value = "2.510.52"
$ws.Range("E18").Value = "  -6.39%  "

# Row 19
Set-TextCell "D19" "11.10"
$ws.Range("E19").Value = "  -6.63%  "

# Row 20
$ws.Range("E20").Value = "  -5.79%  "

# Row 21
Set-TextCell "D21" "323.36"
$ws.Range("E21").Value = "  -6.24%  "

# Row 22
Set-TextCell "D22" "0.966"
$ws.Range("E22").Value = "  -3.38%  "

# Row 23
Set-TextCell "D23" "5.70"
$ws.Range("E23").Value = "  -9.45%  "

# Row 24
Set-TextCell "D24" "60.63"
$ws.Range("E24").Value = "  -4.03%  "

# Row 25
$ws.Range("E25").Value = "  -11.26%  "

# Row 26
$ws.Range("E26").Value = "  -4.92%  "

# Row 27
$ws.Range("E27").Value = "  -2.38%  "

# Row 28
Set-TextCell "D28" "7.65"
$ws.Range("E28").Value = "  -6.67%  "

# Row 29
$ws.Range("E29").Value = "  -6.03%  "

# Row 30
$ws.Range("B30").Value = "PEPE"
$ws.Range("C30").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D30").Value = "0.0₃0768"
$ws.Range("E30").Value = "  -10.44%  "

# Row 31
$ws.Range("B31").Value = "Fetch.AI"
$ws.Range("C31").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
Set-TextCell "D31" "1.25"
$ws.Range("E31").Value = "  -12.92%  "

# Row 32
Set-TextCell "D32" "6.64"
$ws.Range("E32").Value = "  -8.37%  "

# Row 33
Set-TextCell "D33" "0.998"
$ws.Range("E33").Value = "  -0.07%  "

# Row 34
Set-TextCell "D34" "156.22"
$ws.Range("E34").Value = "  -4.82%  "

# Row 35
$ws.Range("B35").Value = "EthereumClassic"
$ws.Range("C35").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextCell "D35" "18.42"
$ws.Range("E35").Value = "  -5.53%  "

# Row 36
$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextCell "D36" "1.37"
$ws.Range("E36").Value = "  -7.94%  "

# Row 37
Set-TextCell "D37" "4.45"
$ws.Range("E37").Value = "  -9.92%  "

# Row 38
$ws.Range("E38").Value = "  -4.49%  "

# Row 39
$ws.Range("B39").Value = "RenderToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextCell "D39" "5.83"
$ws.Range("E39").Value = "  -6.53%  "

# Row 40
$ws.Range("B40").Value = "Bittensor"
$ws.Range("C40").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
Set-TextCell "D40" "313.41"
$ws.Range("E40").Value = "  -10.65%  "

# Row 41
$ws.Range("B41").Value = "OKB"
$ws.Range("C41").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextCell "D41" "36.22"
$ws.Range("E41").Value = "  -5.67%  "

# Row 42
$ws.Range("B42").Value = "SuiNetwork"
$ws.Range("C42").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
Set-TextCell "D42" "0.832"
$ws.Range("E42").Value = "  -12.36%  "

# Row 43
Set-TextCell "D43" "3.70"
$ws.Range("E43").Value = "  -7.34%  "

# Row 44
$ws.Range("E44").Value = "  -0.20%  "

# Row 45
$ws.Range("E45").Value = "  -2.66%  "

# Row 46
Set-TextCell "D46" "0.0941"
$ws.Range("E46").Value = "  -3.06%  "

# Row 47
$ws.Range("E47").Value = "  -6.00%  "

# Row 48
$ws.Range("E48").Value = "  -6.33%  "

# Row 49
Set-TextCell "D49" "0.0229"
$ws.Range("E49").Value = "  -5.50%  "

# Row 50
Set-TextCell "D50" "121.69"
$ws.Range("E50").Value = "  -5.64%  "

# Row 51
Set-TextCell "D51" "18.83"
$ws.Range("E51").Value = "  -10.18%  "
